# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" (column E) / "Valor Mora" (column F) table on rows 16-40
# previously listed periods newest-first (2110 down to 1910). This update
# re-sorts the table oldest-first (1910 up to 2110), carrying each period's
# Valor Mora value along with it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$periods = @(
    "2110","2109","2108","2107","2106","2105","2104","2103","2102","2101",
    "2012","2011","2010","2009","2008","2007","2006","2005","2004","2003",
    "2002","2001","1912","1911","1910"
)

$values = @(
    27604,33125,33125,33125,33125,33125,33125,33125,33125,33125,
    33125,33125,33125,33125,33125,33125,33125,33125,33125,33125,
    33125,33125,33125,33125,17667
)

$firstRow = 16
$lastRow = 40
$count = $lastRow - $firstRow + 1

for ($i = 0; $i -lt $count; $i++) {
    $row = $firstRow + $i
    # Reverse order: new row gets the period/value from the opposite end of the list
    $srcIndex = $count - 1 - $i
    $ws.Range("E$row").Value = $periods[$srcIndex]
    $ws.Range("F$row").Value = $values[$srcIndex]
}
